$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3233.3333
$ws.Range("I64").Value = 3014.2856
$ws.Range("K64").Value = 3014.2856
$ws.Range("M64").Value = -2766.2856
$ws.Range("H67").Value = 3233.3333
$ws.Range("I67").Value = 3014.2856
$ws.Range("K67").Value = 3014.2856
$ws.Range("M67").Value = -2156.2856
$ws.Range("H98").Value = 1465.4615
$ws.Range("I98").Value = 1465.4615
$ws.Range("K98").Value = 1465.4615
$ws.Range("M98").Value = 32.53850000000011
$ws.Range("H112").Value = 1590.9464
$ws.Range("J112").Value = 1605.3273
$ws.Range("L112").Value = 4815.9819
$ws.Range("N112").Value = -7031.9819
$ws.Range("H116").Value = 27025
$ws.Range("I116").Value = 35000
$ws.Range("K116").Value = 35000
$ws.Range("M116").Value = -31558
$ws.Range("H122").Value = 1465.4615
$ws.Range("I122").Value = 1465.4615
$ws.Range("K122").Value = 4396.3845
$ws.Range("M122").Value = -1946.3845
$ws.Range("H131").Value = 1447.3334
$ws.Range("J131").Value = 3800
$ws.Range("L131").Value = 11400
$ws.Range("N131").Value = -21480
$ws.Range("H132").Value = 887.5349
$ws.Range("I132").Value = 814.9729599999999
$ws.Range("J132").Value = 1335
$ws.Range("K132").Value = 2444.91888
$ws.Range("L132").Value = 4005
$ws.Range("M132").Value = 85.08112000000028
$ws.Range("N132").Value = -9065
$ws.Range("H137").Value = 46599.273
$ws.Range("I137").Value = 774.75
$ws.Range("J137").Value = 101588.7
$ws.Range("K137").Value = 2324.25
$ws.Range("L137").Value = 304766.1
$ws.Range("M137").Value = 225.75
$ws.Range("N137").Value = -309866.1
$ws.Range("H138").Value = 1500.86
$ws.Range("I138").Value = 1170.1578
$ws.Range("J138").Value = 1703.5483
$ws.Range("K138").Value = 3510.4734
$ws.Range("L138").Value = 5110.644899999999
$ws.Range("M138").Value = 1629.5266
$ws.Range("N138").Value = -15390.6449
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5629.5107
$ws.Range("I32").Value = 3440.925
$ws.Range("K32").Value = 3440.925
$ws.Range("M32").Value = -3153.925
$ws.Range("H61").Value = 40378.523
$ws.Range("I61").Value = 47837.94
$ws.Range("J61").Value = 8676
$ws.Range("K61").Value = 47837.94
$ws.Range("L61").Value = 8676
$ws.Range("M61").Value = -47625.94
$ws.Range("N61").Value = -9100
$ws.Range("H74").Value = 726.30554
$ws.Range("I74").Value = 734.2
$ws.Range("J74").Value = 450
$ws.Range("K74").Value = 734.2
$ws.Range("L74").Value = 450
$ws.Range("M74").Value = 139.8
$ws.Range("N74").Value = -2198
$ws.Range("H77").Value = 726.30554
$ws.Range("I77").Value = 734.2
$ws.Range("J77").Value = 450
$ws.Range("K77").Value = 3671
$ws.Range("L77").Value = 2250
$ws.Range("M77").Value = 697
$ws.Range("N77").Value = -10986
$ws.Range("H132").Value = 2187.1162
$ws.Range("I132").Value = 2003.579
$ws.Range("K132").Value = 6010.737
$ws.Range("M132").Value = -3480.737
$ws.Range("H136").Value = 40378.523
$ws.Range("I136").Value = 47837.94
$ws.Range("J136").Value = 8676
$ws.Range("K136").Value = 143513.82
$ws.Range("L136").Value = 26028
$ws.Range("M136").Value = -140963.82
$ws.Range("N136").Value = -31128
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5266087.5
$ws.Range("I20").Value = 8335806
$ws.Range("J20").Value = 3713.1428
$ws.Range("K20").Value = 8335806
$ws.Range("L20").Value = 3713.1428
$ws.Range("M20").Value = -8335559
$ws.Range("N20").Value = -4207.1428
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2598.842
$ws.Range("I31").Value = 1396.909
$ws.Range("K31").Value = 1396.909
$ws.Range("M31").Value = -1101.909
$ws.Range("H34").Value = 2598.842
$ws.Range("I34").Value = 1396.909
$ws.Range("K34").Value = 1396.909
$ws.Range("M34").Value = -1194.909
$ws.Range("H99").Value = 2993.3333
$ws.Range("I99").Value = 2990
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 2990
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -1492
$ws.Range("N99").Value = -5996
$ws.Range("H105").Value = 571.75
$ws.Range("I105").Value = 571.75
$ws.Range("K105").Value = 571.75
$ws.Range("M105").Value = 1175.25
$ws.Range("H107").Value = 884.8077
$ws.Range("I107").Value = 884.8077
$ws.Range("K107").Value = 884.8077
$ws.Range("M107").Value = 1035.1923
$ws.Range("H126").Value = 2993.3333
$ws.Range("I126").Value = 2990
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 8970
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -6500
$ws.Range("N126").Value = -13940
$ws.Range("H135").Value = 34694
$ws.Range("J135").Value = 34694
$ws.Range("L135").Value = 34694
$ws.Range("N135").Value = -44834
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 573.0741
$ws.Range("I5").Value = 500
$ws.Range("J5").Value = 894.6
$ws.Range("K5").Value = 1500
$ws.Range("L5").Value = 2683.8
$ws.Range("M5").Value = -1388
$ws.Range("N5").Value = -2907.8
$ws.Range("H107").Value = 509.15384
$ws.Range("J107").Value = 615.5
$ws.Range("L107").Value = 1846.5
$ws.Range("N107").Value = -5686.5
$ws.Range("H122").Value = 1198.6
$ws.Range("J122").Value = 1248.5714
$ws.Range("L122").Value = 11237.1426
$ws.Range("N122").Value = -16137.1426
$ws.Range("H127").Value = 1216.1666
$ws.Range("J127").Value = 1216.1666
$ws.Range("L127").Value = 3648.4998
$ws.Range("N127").Value = -13568.4998
$ws.Range("H131").Value = 42495.895
$ws.Range("J131").Value = 67081.586
$ws.Range("L131").Value = 201244.758
$ws.Range("N131").Value = -211324.758
$ws.Range("H135").Value = 573.0741
$ws.Range("I135").Value = 500
$ws.Range("J135").Value = 894.6
$ws.Range("K135").Value = 4500
$ws.Range("L135").Value = 8051.400000000001
$ws.Range("M135").Value = -1965
$ws.Range("N135").Value = -13121.4
$ws.Range("H141").Value = 3088.5557
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = $null
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3627.75
$ws.Range("I102").Value = 3627.75
$ws.Range("K102").Value = 3627.75
$ws.Range("M102").Value = -2005.75
$ws.Range("H113").Value = 1436.5
$ws.Range("I113").Value = 1436.5
$ws.Range("K113").Value = 1436.5
$ws.Range("M113").Value = 733.5
$ws.Range("H132").Value = 1930460.2
$ws.Range("I132").Value = 3506749.5
$ws.Range("J132").Value = 3884.4443
$ws.Range("K132").Value = 10520248.5
$ws.Range("L132").Value = 11653.3329
$ws.Range("M132").Value = -10517718.5
$ws.Range("N132").Value = -16713.3329
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3347.25
$ws.Range("I7").Value = 3701
$ws.Range("K7").Value = 3701
$ws.Range("M7").Value = -3589
$ws.Range("H40").Value = 9640.091
$ws.Range("I40").Value = 9254.1
$ws.Range("K40").Value = 9254.1
$ws.Range("M40").Value = -9118.1
$ws.Range("H55").Value = 5000547.5
$ws.Range("I55").Value = 9091414
$ws.Range("K55").Value = 9091414
$ws.Range("M55").Value = -9091241
$ws.Range("H61").Value = 2887.5
$ws.Range("I61").Value = 2652.7778
$ws.Range("K61").Value = 2652.7778
$ws.Range("M61").Value = -2450.7778
$ws.Range("H93").Value = 1721.625
$ws.Range("I93").Value = 938.7
$ws.Range("J93").Value = 3026.5
$ws.Range("K93").Value = 938.7
$ws.Range("L93").Value = 3026.5
$ws.Range("M93").Value = 309.3
$ws.Range("N93").Value = -5522.5
$ws.Range("H113").Value = 2887.5
$ws.Range("I113").Value = 2652.7778
$ws.Range("K113").Value = 2652.7778
$ws.Range("M113").Value = -482.7777999999998
$ws.Range("H122").Value = 2718.6365
$ws.Range("I122").Value = 2550.75
$ws.Range("J122").Value = 3166.3333
$ws.Range("K122").Value = 7652.25
$ws.Range("L122").Value = 9498.999899999999
$ws.Range("M122").Value = -5202.25
$ws.Range("N122").Value = -14398.9999
$ws.Range("H126").Value = 3347.25
$ws.Range("I126").Value = 3701
$ws.Range("K126").Value = 11103
$ws.Range("M126").Value = -8633
$ws.Range("H132").Value = 4301.4736
$ws.Range("I132").Value = 1984.1428
$ws.Range("K132").Value = 5952.428400000001
$ws.Range("M132").Value = -3422.428400000001
$ws.Range("H136").Value = 3598.7778
$ws.Range("I136").Value = 2631.8333
$ws.Range("J136").Value = 5532.6665
$ws.Range("K136").Value = 7895.499899999999
$ws.Range("L136").Value = 16597.9995
$ws.Range("M136").Value = -5345.499899999999
$ws.Range("N136").Value = -21697.9995
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1520.4
$ws.Range("I113").Value = 1349.5
$ws.Range("J113").Value = 1634.3334
$ws.Range("K113").Value = 4048.5
$ws.Range("L113").Value = 4903.0002
$ws.Range("M113").Value = -1878.5
$ws.Range("N113").Value = -9243.0002
$ws.Range("H132").Value = 1189.9181
$ws.Range("I132").Value = 1083.5306
$ws.Range("J132").Value = 1624.3334
$ws.Range("K132").Value = 3250.5918
$ws.Range("L132").Value = 4873.0002
$ws.Range("M132").Value = -720.5918000000001
$ws.Range("N132").Value = -9933.0002
$ws.Range("H136").Value = 15874222
$ws.Range("I136").Value = 22223002
$ws.Range("J136").Value = 2269.4
$ws.Range("K136").Value = 66669006
$ws.Range("L136").Value = 6808.200000000001
$ws.Range("M136").Value = -66666456
$ws.Range("N136").Value = -11908.2
$ws.Range("H141").Value = 68695.06
$ws.Range("J141").Value = 69608.07000000001
$ws.Range("L141").Value = 69608.07000000001
$ws.Range("N141").Value = -79968.07000000001
